# Updated cryptos list on Mon Nov 25 10:18:20 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for the crypto list, and
# swaps the PancakeSwap/Hedera (rows 30-31) and Kaspa/Bittensor
# (rows 39-40) rows whose ranking order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.614.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.473.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.51'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '669.65'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.50'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.432'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.470.93'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.28'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +11.05%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '98.434.14'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.20'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000260'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.134.16'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.478.67'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.18%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.538'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.72'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +6.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '521.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.45'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000204'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.86'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +5.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '94.68'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.68'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.56'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +9.04%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.146'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.88'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +10.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.191'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.584'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '30.24'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.07'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.53'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.156'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '531.90'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.892'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.79'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0438'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.43'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.76'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.67'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.26'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +8.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.78'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.26'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.84%  '
